# Refresh the flight search results grid (columns A:D, rows 2-19) with the
# latest scraped values (airline, departure time, duration, price).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SpiceJet",  "21:40", "02 h 45 m", "₹ 8,982"),
    @("SpiceJet",  "08:30", "05 h 30 m", "₹ 8,982"),
    @("Air India", "06:00", "02 h 55 m", "₹ 8,983"),
    @("IndiGo",    "06:20", "02 h 50 m", "₹ 8,983"),
    @("Vistara",   "07:05", "02 h 50 m", "₹ 8,983"),
    @("IndiGo",    "08:45", "02 h 45 m", "₹ 8,983"),
    @("Air India", "09:55", "03 h 05 m", "₹ 8,983"),
    @("Vistara",   "10:35", "02 h 50 m", "₹ 8,983"),
    @("IndiGo",    "10:40", "02 h 45 m", "₹ 8,983"),
    @("IndiGo",    "13:20", "03 h",      "₹ 8,983"),
    @("IndiGo",    "15:10", "02 h 50 m", "₹ 8,983"),
    @("IndiGo",    "16:35", "02 h 50 m", "₹ 8,983"),
    @("Air India", "16:55", "02 h 55 m", "₹ 8,983"),
    @("Vistara",   "17:15", "02 h 45 m", "₹ 8,983"),
    @("IndiGo",    "18:15", "02 h 50 m", "₹ 8,983"),
    @("IndiGo",    "19:35", "02 h 50 m", "₹ 8,983"),
    @("Air India", "20:15", "02 h 30 m", "₹ 8,983"),
    @("Vistara",   "21:05", "02 h 50 m", "₹ 8,983")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
